$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns updated per row, in order: B, C, D, F, G, H, I, K, M
$columns = @("B", "C", "D", "F", "G", "H", "I", "K", "M")

$data = @{
    2 = @(0.1787717633110617, 0.03353462622342818, 0.02892823157375446, 0.9168275556501868, 0.7678336403582477, 0.8280886392317086, 0.7021028952361803, 0.1917756758864329, 0.8825804115905527)
    3 = @(0.1575630821617438, 0.02956669347233287, 0.02807602933530262, 0.9058254434294426, 0.757775792964253, 0.8275313332253234, 0.700236497951721, 0.1676882481611557, 0.7860361702641683)
    4 = @(0.1445492513491473, 0.02711624283243452, 0.02754214540397371, 0.8995832456350428, 0.7520574427439897, 0.827543678885192, 0.6994339721324465, 0.1528842266068438, 0.727296995630752)
    5 = @(0.1392484981522415, 0.02611414034200266, 0.02732193677843142, 0.8971682413484103, 0.7498417245989515, 0.8276377711995764, 0.6991930878981591, 0.1468482556740724, 0.7034896055670004)
    6 = @(0.1383684736421884, 0.0259475299682066, 0.02728521222221048, 0.896774998933104, 0.749480713158917, 0.827658771668041, 0.6991582864298778, 0.1458458042228727, 0.6995440331435958)
    7 = @(0.1444777529572434, 0.02710274236048349, 0.02753918627547591, 0.8995501551519993, 0.7520270975021788, 0.8275445873543674, 0.6994303749352397, 0.1528028359060727, 0.7269754055958657)
    8 = @(0.1714574905424513, 0.03216942981165971, 0.02863661347872437, 0.9129273425157152, 0.764270547080045, 0.8278228407560704, 0.7013879486869854, 0.1834735010080664, 0.849176137151332)
    9 = @(0.2244171656265053, 0.04199236852471699, 0.03070327271258222, 0.9432472336838487, 0.7919284706224659, 0.8311865309906779, 0.7079624058502247, 0.2434922760148481, 1.093392505533117)
    10 = @(0.2633436677054988, 0.04914016817168942, 0.03216827989351856, 0.9680403655024321, 0.8145049845615517, 0.835384050270477, 0.7144764332089082, 0.2874975487356153, 1.276051646617873)
    11 = @(0.2810531195605677, 0.05237685049391416, 0.03282289162630292, 0.9798717050848751, 0.8252726701581139, 0.8376703418863087, 0.7178090246420581, 0.3074944418446819, 1.359948983987636)
    12 = @(0.2877591269201218, 0.05360033860029034, 0.03306905029229057, 0.9844318158732079, 0.8294221732038949, 0.8385904202993117, 0.719124366054487, 0.3150633642891876, 1.391842403467592)
    13 = @(0.2863148826108954, 0.05333693571694198, 0.03301611294048001, 0.9834461580277889, 0.8285252930686227, 0.8383898475389344, 0.7188387064374666, 0.3134334219727748, 1.384967977314147)
    14 = @(0.2816048327064209, 0.05247755133648013, 0.03284317805050563, 0.9802452660613881, 0.8256126062823057, 0.8377449481367734, 0.7179161676760586, 0.3081172134220367, 1.362570359115239)
    15 = @(0.2787197559569279, 0.05195086960605977, 0.03273702449931903, 0.9782950350169415, 0.8238378939713868, 0.8373570050115262, 0.7173580428052944, 0.3048604188167587, 1.348867459697644)
    16 = @(0.2621863106109004, 0.04892834102184906, 0.03212525909869157, 0.9672783131500893, 0.8138113425397364, 0.8352422287904631, 0.7142660938675149, 0.2861902413606572, 1.270585597070422)
    17 = @(0.2520436951929526, 0.0470702814920827, 0.03174691172713295, 0.9606617455226996, 0.8077881384325707, 0.8340414815359622, 0.7124640559212096, 0.2747309448907345, 1.222773673458235)
    18 = @(0.2462101050033425, 0.04600017544144919, 0.03152818473250818, 0.9569080840421265, 0.8043705508832772, 0.8333863039837439, 0.7114623190031821, 0.2681378771782192, 1.195348803229635)
    19 = @(0.2442349938823156, 0.0456376165667649, 0.03145393746231662, 0.9556460816328496, 0.8032214381067888, 0.8331705587702771, 0.7111291083766673, 0.2659052492898013, 1.186075915304258)
    20 = @(0.2531233788656664, 0.04726822016293397, 0.03178730269826957, 0.9613607047881345, 0.8084244717495892, 0.8341656322177897, 0.712652287867293, 0.2759510147465107, 1.227855488399342)
    21 = @(0.282988296591725, 0.05273003259999598, 0.0328940203541066, 0.9811832760455133, 0.8264661747700757, 0.8379328957892511, 0.7181856894361687, 0.3096788100183687, 1.369145674187223)
    22 = @(0.3025055232419902, 0.05628696039076431, 0.03360723593487336, 0.9946039989087723, 0.8386774463912019, 0.8407116357415134, 0.7221132370302783, 0.3317014806761165, 1.462209439755611)
    23 = @(0.2920890579669901, 0.05438973224138977, 0.0332275116111731, 0.9873983960979871, 0.8321214761160149, 0.839199558317091, 0.7199884766644544, 0.3199495648623554, 1.412470890922805)
    24 = @(0.2526352615255689, 0.04717873795311789, 0.03176904570707251, 0.9610445486837733, 0.8081366446028255, 0.8341093942252513, 0.7125670814530167, 0.2753994366006509, 1.225557804698937)
    25 = @(0.2100861053138772, 0.03934711647532652, 0.0301534644506134, 0.9346047141559524, 0.7840523722920523, 0.8299740256247787, 0.7058893187736572, 0.2272705328956732, 1.026790146739458)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $cellRef = "$($columns[$i])$row"
        $ws.Range($cellRef).Value = $values[$i]
    }
}
